$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared-string text fixups (masthead volume/issue number and date range) ---
$ws.Range("A8").Value = "Volume 32   Number  11"
$ws.Range("C9").Value = "Report Covering the Week  3/10/2025  Through  3/16/2025"

# --- Simple value updates (number format / style unchanged) ---
$ws.Range("I15").Value = 5
$ws.Range("K15").Value = 400
$ws.Range("L15").Value = 150
$ws.Range("M15").Value = 400
$ws.Range("N15").Value = -28.571428571428
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 9
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = -10
$ws.Range("I16").Value = 25
$ws.Range("J16").Value = 38
$ws.Range("K16").Value = -34.210526315789
$ws.Range("L16").Value = 8.695652173913
$ws.Range("M16").Value = 8.695652173913
$ws.Range("N16").Value = -81.343283582089
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = -50
$ws.Range("F17").Value = 13
$ws.Range("G17").Value = 11
$ws.Range("H17").Value = 18.181818181818
$ws.Range("I17").Value = 42
$ws.Range("J17").Value = 32
$ws.Range("K17").Value = 31.25
$ws.Range("L17").Value = 20
$ws.Range("M17").Value = 55.555555555555
$ws.Range("N17").Value = -10.63829787234
$ws.Range("C18").Value = 2
$ws.Range("E18").Value = -33.333333333333
$ws.Range("G18").Value = 14
$ws.Range("H18").Value = -42.857142857142
$ws.Range("I18").Value = 26
$ws.Range("J18").Value = 25
$ws.Range("K18").Value = 4
$ws.Range("L18").Value = -18.75
$ws.Range("M18").Value = -29.729729729729
$ws.Range("N18").Value = -88.646288209607
$ws.Range("C19").Value = 9
$ws.Range("E19").Value = -35.714285714285
$ws.Range("F19").Value = 33
$ws.Range("G19").Value = 57
$ws.Range("H19").Value = -42.105263157894
$ws.Range("I19").Value = 91
$ws.Range("J19").Value = 160
$ws.Range("K19").Value = -43.125
$ws.Range("L19").Value = -24.166666666666
$ws.Range("M19").Value = 56.896551724137
$ws.Range("N19").Value = 30
$ws.Range("C20").Value = 9
$ws.Range("D20").Value = 7
$ws.Range("E20").Value = 28.571428571428
$ws.Range("F20").Value = 19
$ws.Range("G20").Value = 31
$ws.Range("H20").Value = -38.709677419354
$ws.Range("I20").Value = 47
$ws.Range("J20").Value = 86
$ws.Range("K20").Value = -45.348837209302
$ws.Range("L20").Value = -32.857142857142
$ws.Range("M20").Value = 161.111111111111
$ws.Range("N20").Value = -85.757575757575
$ws.Range("C21").Value = 26
$ws.Range("D21").Value = 32
$ws.Range("E21").Value = -18.75
$ws.Range("F21").Value = 83
$ws.Range("G21").Value = 123
$ws.Range("H21").Value = -32.520325203252
$ws.Range("I21").Value = 236
$ws.Range("J21").Value = 343
$ws.Range("K21").Value = -31.195335276967
$ws.Range("L21").Value = -16.312056737588
$ws.Range("M21").Value = 41.317365269461
$ws.Range("N21").Value = -71.184371184371
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 1
$ws.Range("G23").Value = 5
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 13
$ws.Range("J23").Value = 16
$ws.Range("K23").Value = -18.75
$ws.Range("L23").Value = -18.75
$ws.Range("M23").Value = 30
$ws.Range("C24").Value = 19
$ws.Range("D24").Value = 16
$ws.Range("E24").Value = 18.75
$ws.Range("F24").Value = 65
$ws.Range("G24").Value = 78
$ws.Range("H24").Value = -16.666666666666
$ws.Range("I24").Value = 178
$ws.Range("J24").Value = 215
$ws.Range("K24").Value = -17.209302325581
$ws.Range("L24").Value = -12.31527093596
$ws.Range("M24").Value = 17.105263157894
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = -100
$ws.Range("F25").Value = 10
$ws.Range("G25").Value = 25
$ws.Range("H25").Value = -60
$ws.Range("J25").Value = 79
$ws.Range("K25").Value = -50.632911392405
$ws.Range("L25").Value = -53.012048192771
$ws.Range("C26").Value = 3
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = -57.142857142857
$ws.Range("F26").Value = 17
$ws.Range("G26").Value = 23
$ws.Range("H26").Value = -26.086956521739
$ws.Range("I26").Value = 64
$ws.Range("J26").Value = 63
$ws.Range("K26").Value = 1.587301587301
$ws.Range("L26").Value = 14.285714285714
$ws.Range("M26").Value = -5.882352941176
$ws.Range("I27").Value = 6
$ws.Range("K27").Value = 50
$ws.Range("L27").Value = 20
$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 7
$ws.Range("H28").Value = -42.857142857142
$ws.Range("I28").Value = 7
$ws.Range("J28").Value = 11
$ws.Range("K28").Value = -36.363636363636
$ws.Range("L28").Value = -12.5

# --- Updates that also change the cell style (text <-> number) ---
# Strategy: write the new literal via Formula (so text stays text), then
# PasteSpecial formats-only from a stable donor cell that already carries
# the desired style, avoiding creation of brand-new style entries.
$ws.Range("C15").Formula = "1"
$ws.Range("J14").Copy()
$ws.Range("C15").PasteSpecial(-4122)

$ws.Range("F15").Formula = "1"
$ws.Range("J14").Copy()
$ws.Range("F15").PasteSpecial(-4122)

$ws.Range("C25").Formula = "'0"
$ws.Range("D30").Copy()
$ws.Range("C25").PasteSpecial(-4122)

$ws.Range("C27").Formula = "1"
$ws.Range("J14").Copy()
$ws.Range("C27").PasteSpecial(-4122)

$ws.Range("F27").Formula = "1"
$ws.Range("J14").Copy()
$ws.Range("F27").PasteSpecial(-4122)

$ws.Range("C28").Formula = "2"
$ws.Range("J14").Copy()
$ws.Range("C28").PasteSpecial(-4122)

$ws.Range("D28").Formula = "4"
$ws.Range("J14").Copy()
$ws.Range("D28").PasteSpecial(-4122)

$ws.Range("E28").Formula = "-50"
$ws.Range("K14").Copy()
$ws.Range("E28").PasteSpecial(-4122)

$ws.Range("D31").Formula = "1"
$ws.Range("J14").Copy()
$ws.Range("D31").PasteSpecial(-4122)

$ws.Range("E31").Formula = "-100"
$ws.Range("K14").Copy()
$ws.Range("E31").PasteSpecial(-4122)

$ws.Range("G31").Formula = "1"
$ws.Range("J14").Copy()
$ws.Range("G31").PasteSpecial(-4122)

$ws.Range("H31").Formula = "-100"
$ws.Range("K14").Copy()
$ws.Range("H31").PasteSpecial(-4122)

$ws.Range("J31").Formula = "1"
$ws.Range("J14").Copy()
$ws.Range("J31").PasteSpecial(-4122)

$ws.Range("K31").Formula = "-100"
$ws.Range("K14").Copy()
$ws.Range("K31").PasteSpecial(-4122)

